$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7,3,3,0),
    @(2,3,3,1),
    @(4,0,3,2),
    @(7,3,7,0),
    @(5,1,6,2),
    @(3,1,3,2),
    @(4,1,5,2),
    @(2,2,3,1),
    @(4,0,5,3),
    @(3,1,2,2),
    @(4,0,4,3),
    @(6,1,5,2),
    @(5,0,6,2),
    @(4,2,4,0),
    @(4,0,5,3),
    @(3,0,5,3),
    @(4,0,6,2),
    @(7,2,5,1),
    @(5,2,6,1),
    @(2,2,4,0),
    @(5,1,4,2),
    @(2,2,4,1)
)

$startRow = 1457
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$lastRow = $r - 1
$nextRow = $lastRow + 1

$ws.Range("A$nextRow").Select() | Out-Null
